# Update "按行业分连锁零售企业年末从业人数" sheet:
#  - Drop the oldest five years of data (2005年-2009年), shifting the
#    2010年-2020年 rows up to rows 2-12.
#  - Append a new 2021年 row (row 13) with its own figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2005年-2009年 rows (original rows 2-6); remaining rows shift up.
$ws.Range("A2:A6").EntireRow.Delete()

# After the shift, 2020年 is on row 12. Clone its formatting (style) down to
# row 13 so the new year label keeps the same look (border/alignment/bold).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Fill in the new 2021年 figures on row 13.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 0.8169
$ws.Cells.Item(13, 3).Value = 47.936
$ws.Cells.Item(13, 4).Value = 6.8172
$ws.Cells.Item(13, 5).Value = 4.554
$ws.Cells.Item(13, 6).Value = 0.1411
$ws.Cells.Item(13, 7).Value = 12.1276
$ws.Cells.Item(13, 8).Value = 13.9024
$ws.Cells.Item(13, 9).Value = 111.8904
$ws.Cells.Item(13, 10).Value = 6.6445
